$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BLACK MYTH IMP
$ws.Range("A2").Value = "BLACK MYTH IMP"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1182662
$ws.Range("E2").Value = 1182662
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 45799.71182745069
$ws.Range("G2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 3: Lukedom
$ws.Range("A3").Value = "Lukedom"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 1549600
$ws.Range("D3").Value = 604308
$ws.Range("E3").Value = 2153908
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 45799.7118275649
$ws.Range("G3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
